$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subject codes and names, plus year on row 5
$ws.Range("C2").Value = "SEP401"
$ws.Range("D2").Value = "Software Engineering Principles"

$ws.Range("C3").Value = "SDM404"
$ws.Range("D3").Value = "Software Development Management"

$ws.Range("C4").Value = "REM502"
$ws.Range("D4").Value = "Research Methodologies"

$ws.Range("A5").Value = 2024.0
$ws.Range("C5").Value = "Elective 1"
$ws.Range("D5").Value = "Elective Subject"

# Apply font (Arial, theme color 1, no explicit size) to columns C and D, rows 2-5
$range = $ws.Range("C2:D5")
$range.Font.Name = "Arial"

# Set column D width (raw OOXML width ~27.63 characters)
$ws.Columns.Item(4).ColumnWidth = 26.83
